$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 6, leaving only the header row and the first data row
$ws.Range("A3:G6").EntireRow.Delete()

# Update row 2 with new data
$ws.Range("A2").Value = "Andrian Putra "
$ws.Range("B2").Value = 21212344
$ws.Range("C2").Value = "ramadanand89@gmail.com"
$ws.Range("D2").Value = 3232
$ws.Range("E2").Value = "VII - A"
$ws.Range("F2").Value = "eItFafz5lj"
$ws.Range("G2").Value = 12345678
